# Apply "Updated symbol list" edits to cryptos.xlsx (sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''248.46'

$ws.Range("D3").Value = '''22.67'

$ws.Range("D4").Value = '''5.277'

$ws.Range("D5").Value = '''0.05684'

$ws.Range("D6").Value = '''3.420'

$ws.Range("D7").Value = '''6.328'

$ws.Range("D8").Value = '''0.8074'

$ws.Range("D9").Value = '''0.8946'

$ws.Range("D11").Value = '''0.07476'

$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = '''0.03054'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'

$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.03098'
$ws.Range("E13").Value = '12BitrueCoinBTR'

$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09396'
$ws.Range("E14").Value = '13BitMartTokenBMX'

$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = '''3.881'
$ws.Range("E15").Value = '14MCDexMCB'

$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = '''0.001581'
$ws.Range("E16").Value = '15BitForexTokenBF'

$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = '''0.04788'
$ws.Range("E17").Value = '16CoinExTokenCET'

$ws.Range("B18").Value = 'UpBots'
$ws.Range("C18").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D18").Value = '''0.01830'
$ws.Range("E18").Value = '17UpBotsUBXTBestin24h'

$ws.Range("B19").Value = 'One'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D19").Value = '''0.0005806'
$ws.Range("E19").Value = '18OneONE'

$ws.Range("B20").Value = 'TigerCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D20").Value = '''0.006418'
$ws.Range("E20").Value = '19TigerCashTCH'

$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D21").Value = '''0.004967'
$ws.Range("E21").Value = '20HotbitTokenHTB'

$ws.Range("B22").Value = 'BitKan'
$ws.Range("C22").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D22").Value = '''0.0009972'
$ws.Range("E22").Value = '21BitKanKAN'

$ws.Range("B23").Value = 'NitroEx'
$ws.Range("C23").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D23").Value = '''0.0001501'
$ws.Range("E23").Value = '22NitroExNTX'

$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D24").Value = '''3.685'
$ws.Range("E24").Value = '23LEOLEO'

$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").Value = '''2.158'
$ws.Range("E25").Value = '24BTSETokenBTSE'

$ws.Range("B26").Value = 'BitpandaEcosystemToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D26").Value = '''0.3256'
$ws.Range("E26").Value = '25BitpandaEcosystemTokenBEST'

$ws.Range("B27").Value = 'ProBitToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D27").Value = '''0.1332'
$ws.Range("E27").Value = '26ProBitTokenPROB'

$ws.Range("D40").Value = '''0.03958'

$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '''0.1070'
$ws.Range("E41").Value = '40BKEXTokenBKK'

$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = '''0.002731'
$ws.Range("E42").Value = '41CEJICEJI'

$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = '''0.003041'
$ws.Range("E43").Value = '42KickTokenKICKWorstin24h'

$ws.Range("D44").Value = '''0.008763'

$ws.Range("D45").Value = '''0.00005581'

$ws.Range("D47").Value = '''0.4995'

$ws.Range("D48").Value = '''0.2026'
$ws.Range("E48").Value = '47BOLOBOLO'

Write-Output "Applied cryptos.xlsx symbol/price updates"
